$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("HOTCARD")

$ws.Range("B21").Value = "First National Bank"
$ws.Range("B21").Style = "Normal"

$ws.Range("D21").Value = "FISB"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "Legacy Core"
$ws.Range("E21").Style = "Normal"

$ws.Range("F21").Value = "PaymentsOne Debit"
$ws.Range("F21").Style = "Normal"

$ws.Range("H21").Value = "Dallas, TX"
$ws.Range("H21").Style = "Normal"

$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "15000"
$ws.Range("I21").Style = "Normal"

$ws.Range("J21").Value = "Yes"
$ws.Range("J21").Style = "Normal"
